# Applies the "Make Mail Service Functionality" edit:
#  1. Rename the sheet from "07.09.2020" to "05.08.2020" (this also updates
#     the chart's series/category formula references since they are
#     sheet-name-qualified).
#  2. Populate row 6 (Prognoza energie [MWh]) with forecast values.
#  3. Populate row 7 (Consum realizat [MWh]) with actual consumption values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet -- Excel automatically re-points any formulas
#    (including chart series references) that referred to the old name.
$ws.Name = "05.08.2020"

# 2. Row 6 values (C6:Z6) -- all 2.5
$row6Values = @(2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5)

# 3. Row 7 values (C7:Z7)
$row7Values = @(2.312,2.296,2.295,2.29,2.287,2.248,2.341,2.383,1.534,0.918,2.171,2.42,2.438,2.393,2.344,2.442,2.479,2.41,2.353,2.45,2.499,2.31,2.244,2.238)

$startCol = 3   # column C

for ($i = 0; $i -lt $row6Values.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(6, $col).Value = $row6Values[$i]
}

for ($i = 0; $i -lt $row7Values.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(7, $col).Value = $row7Values[$i]
}

# 4. The chart's series formulas are hard-coded with the sheet name, so make
#    sure they are re-pointed at the newly renamed sheet as well.
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart

$series1 = $chart.SeriesCollection().Item(1)
$series1.Formula = '=SERIES("Prognoza energie [MWh]",' + "'05.08.2020'" + '!$C$5:$Z$5,' + "'05.08.2020'" + '!$C$6:$Z$6,1)'

$series2 = $chart.SeriesCollection().Item(2)
$series2.Formula = '=SERIES("Consum realizat [MWh]",' + "'05.08.2020'" + '!$C$5:$Z$5,' + "'05.08.2020'" + '!$C$7:$Z$7,2)'
